$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 26.81310566666667
$ws.Range("H2").Value = 80.439317
$ws.Range("I2").Value = 0.004518206005002021
$ws.Range("J2").Value = 0.004518206005002021
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2457116666666667
$ws.Range("N2").Value = 0.737135
$ws.Range("O2").Value = 0.007550096805344261
$ws.Range("P2").Value = 0.007550096805344259
$ws.Range("Q2").Value = 6.588292881866112
$ws.Range("R2").Value = 59.294635936795
$ws.Range("S2").Value = 0.00003411289272425302
$ws.Range("T2").Value = 0.000034112892724253
$ws.Range("G3").Value = 26.81310566666667
$ws.Range("H3").Value = 80.439317
$ws.Range("I3").Value = 0.004518206005002021
$ws.Range("J3").Value = 0.004518206005002021
$ws.Range("N3").Value = 94.553567
$ws.Range("O3").Value = 0.9684638283904637
$ws.Range("P3").Value = 0.9684638283904636
$ws.Range("Q3").Value = 845.0915943770822
$ws.Range("R3").Value = 7605.824349393739
$ws.Range("S3").Value = 0.00437571908506104
$ws.Range("T3").Value = 0.004375719085061039
$ws.Range("G4").Value = 26.81310566666667
$ws.Range("H4").Value = 80.439317
$ws.Range("I4").Value = 0.004518206005002021
$ws.Range("J4").Value = 0.004518206005002021
$ws.Range("M4").Value = 0.737729
$ws.Range("N4").Value = 2.213187
$ws.Range("O4").Value = 0.02266854253064832
$ws.Range("P4").Value = 0.02266854253064832
$ws.Range("Q4").Value = 19.78080563036433
$ws.Range("R4").Value = 178.027250673279
$ws.Range("S4").Value = 0.0001024211449866189
$ws.Range("T4").Value = 0.0001024211449866189
$ws.Range("G5").Value = 26.81310566666667
$ws.Range("H5").Value = 80.439317
$ws.Range("I5").Value = 0.004518206005002021
$ws.Range("J5").Value = 0.004518206005002021
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.042878
$ws.Range("N5").Value = 0.128634
$ws.Range("O5").Value = 0.001317532273543725
$ws.Range("P5").Value = 0.001317532273543725
$ws.Range("Q5").Value = 1.149692344775333
$ws.Range("R5").Value = 10.347231102978
$ws.Range("S5").Value = 0.000005952882230109224
$ws.Range("T5").Value = 0.000005952882230109222
$ws.Range("G6").Value = 5771.873535333333
$ws.Range("I6").Value = 0.9726032482643521
$ws.Range("J6").Value = 0.9726032482643523
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2457116666666667
$ws.Range("N6").Value = 0.737135
$ws.Range("O6").Value = 0.007550096805344261
$ws.Range("P6").Value = 0.007550096805344259
$ws.Range("Q6").Value = 1418.216666155979
$ws.Range("R6").Value = 12763.94999540381
$ws.Range("S6").Value = 0.007343248677588136
$ws.Range("T6").Value = 0.007343248677588136
$ws.Range("G7").Value = 5771.873535333333
$ws.Range("I7").Value = 0.9726032482643521
$ws.Range("J7").Value = 0.9726032482643523
$ws.Range("N7").Value = 94.553567
$ws.Range("O7").Value = 0.9684638283904637
$ws.Range("P7").Value = 0.9684638283904636
$ws.Range("S7").Value = 0.9419310653190951
$ws.Range("T7").Value = 0.9419310653190952
$ws.Range("G8").Value = 5771.873535333333
$ws.Range("I8").Value = 0.9726032482643521
$ws.Range("J8").Value = 0.9726032482643523
$ws.Range("M8").Value = 0.737729
$ws.Range("N8").Value = 2.213187
$ws.Range("O8").Value = 0.02266854253064832
$ws.Range("P8").Value = 0.02266854253064832
$ws.Range("Q8").Value = 4258.078491347925
$ws.Range("R8").Value = 38322.70642213133
$ws.Range("S8").Value = 0.02204749809872717
$ws.Range("T8").Value = 0.02204749809872717
$ws.Range("G9").Value = 5771.873535333333
$ws.Range("I9").Value = 0.9726032482643521
$ws.Range("J9").Value = 0.9726032482643523
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.042878
$ws.Range("N9").Value = 0.128634
$ws.Range("O9").Value = 0.001317532273543725
$ws.Range("P9").Value = 0.001317532273543725
$ws.Range("Q9").Value = 247.4863934480227
$ws.Range("R9").Value = 2227.377541032204
$ws.Range("S9").Value = 0.001281436168941744
$ws.Range("T9").Value = 0.001281436168941744
$ws.Range("G10").Value = 132.4457753333333
$ws.Range("H10").Value = 397.337326
$ws.Range("I10").Value = 0.02231808970163987
$ws.Range("J10").Value = 0.02231808970163988
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.2457116666666667
$ws.Range("N10").Value = 0.737135
$ws.Range("O10").Value = 0.007550096805344261
$ws.Range("P10").Value = 0.007550096805344259
$ws.Range("Q10").Value = 32.54347220011222
$ws.Range("R10").Value = 292.89124980101
$ws.Range("S10").Value = 0.0001685037377577379
$ws.Range("T10").Value = 0.0001685037377577378
$ws.Range("G11").Value = 132.4457753333333
$ws.Range("H11").Value = 397.337326
$ws.Range("I11").Value = 0.02231808970163987
$ws.Range("J11").Value = 0.02231808970163988
$ws.Range("N11").Value = 94.553567
$ws.Range("O11").Value = 0.9684638283904637
$ws.Range("P11").Value = 0.9684638283904636
$ws.Range("Q11").Value = 4174.40683061576
$ws.Range("R11").Value = 37569.66147554184
$ws.Range("S11").Value = 0.02161426259481193
$ws.Range("T11").Value = 0.02161426259481193
$ws.Range("G12").Value = 132.4457753333333
$ws.Range("H12").Value = 397.337326
$ws.Range("I12").Value = 0.02231808970163987
$ws.Range("J12").Value = 0.02231808970163988
$ws.Range("M12").Value = 0.737729
$ws.Range("N12").Value = 2.213187
$ws.Range("O12").Value = 0.02266854253064832
$ws.Range("P12").Value = 0.02266854253064832
$ws.Range("Q12").Value = 97.70908939088466
$ws.Range("R12").Value = 879.381804517962
$ws.Range("S12").Value = 0.0005059185656044477
$ws.Range("T12").Value = 0.0005059185656044478
$ws.Range("G13").Value = 132.4457753333333
$ws.Range("H13").Value = 397.337326
$ws.Range("I13").Value = 0.02231808970163987
$ws.Range("J13").Value = 0.02231808970163988
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.042878
$ws.Range("N13").Value = 0.128634
$ws.Range("O13").Value = 0.001317532273543725
$ws.Range("P13").Value = 0.001317532273543725
$ws.Range("Q13").Value = 5.679009954742666
$ws.Range("R13").Value = 51.111089592684
$ws.Range("S13").Value = 0.00002940480346575438
$ws.Range("T13").Value = 0.00002940480346575438
$ws.Range("G14").Value = 3.326003
$ws.Range("H14").Value = 9.978009
$ws.Range("I14").Value = 0.0005604560290058679
$ws.Range("J14").Value = 0.000560456029005868
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.2457116666666667
$ws.Range("N14").Value = 0.737135
$ws.Range("O14").Value = 0.007550096805344261
$ws.Range("P14").Value = 0.007550096805344259
$ws.Range("Q14").Value = 0.8172377404683333
$ws.Range("R14").Value = 7.355139664215
$ws.Range("S14").Value = 0.000004231497274133133
$ws.Range("T14").Value = 0.000004231497274133133
$ws.Range("G15").Value = 3.326003
$ws.Range("H15").Value = 9.978009
$ws.Range("I15").Value = 0.0005604560290058679
$ws.Range("J15").Value = 0.000560456029005868
$ws.Range("N15").Value = 94.553567
$ws.Range("O15").Value = 0.9684638283904637
$ws.Range("P15").Value = 0.9684638283904636
$ws.Range("Q15").Value = 104.8284825009003
$ws.Range("R15").Value = 943.456342508103
$ws.Range("S15").Value = 0.0005427813914955396
$ws.Range("T15").Value = 0.0005427813914955396
$ws.Range("G16").Value = 3.326003
$ws.Range("H16").Value = 9.978009
$ws.Range("I16").Value = 0.0005604560290058679
$ws.Range("J16").Value = 0.000560456029005868
$ws.Range("M16").Value = 0.737729
$ws.Range("N16").Value = 2.213187
$ws.Range("O16").Value = 0.02266854253064832
$ws.Range("P16").Value = 0.02266854253064832
$ws.Range("Q16").Value = 2.453688867187
$ws.Range("R16").Value = 22.083199804683
$ws.Range("S16").Value = 0.00001270472133007778
$ws.Range("T16").Value = 0.00001270472133007779
$ws.Range("G17").Value = 3.326003
$ws.Range("H17").Value = 9.978009
$ws.Range("I17").Value = 0.0005604560290058679
$ws.Range("J17").Value = 0.000560456029005868
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.042878
$ws.Range("N17").Value = 0.128634
$ws.Range("O17").Value = 0.001317532273543725
$ws.Range("P17").Value = 0.001317532273543725
$ws.Range("Q17").Value = 0.142612356634
$ws.Range("R17").Value = 1.283511209706
$ws.Range("S17").Value = 0.000000738418906117389
$ws.Range("T17").Value = 0.000000738418906117389
